$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.589.90'
$ws.Range("E2").Value = '  -0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.976.36'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.63'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("E6").Value = '  +2.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.14'
$ws.Range("E7").Value = '  +2.78%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.379'
$ws.Range("E9").Value = '  +1.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0790'
$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.26'
$ws.Range("E12").Value = '  +3.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.843'
$ws.Range("E13").Value = '  +2.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.267.88'
$ws.Range("E14").Value = '  +0.46%  '
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.44'
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.981.89'
$ws.Range("E17").Value = '  +0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.606.61'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.88'
$ws.Range("E19").Value = '  +0.29%  '
$ws.Range("E20").Value = '  -0.38%  '
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.91'
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("E24").Value = '  +2.53%  '
$ws.Range("E25").Value = '  +1.62%  '
$ws.Range("E26").Value = '  +7.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.19'
$ws.Range("E27").Value = '  -2.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.95'
$ws.Range("E28").Value = '  +0.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.39'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  +19.43%  '
$ws.Range("E31").Value = '  +1.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.80'
$ws.Range("E32").Value = '  +2.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0618'
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("E34").Value = '  +6.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.28'
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("E36").Value = '  +0.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("E38").Value = '  -2.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.42'
$ws.Range("E39").Value = '  -11.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0973'
$ws.Range("E40").Value = '  -2.53%  '
$ws.Range("E41").Value = '  +0.52%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("E43").Value = '  -0.70%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.97'
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.365.89'
$ws.Range("E45").Value = '  -0.20%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.20'
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("E50").Value = '  +6.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.161.91'
$ws.Range("E51").Value = '  +0.65%  '
